$d = $word.ActiveDocument

# --- Body content updates -------------------------------------------------

# 1. Bump the patch version referenced in the intro sentence:
#    "These are known issues in Lightning 1.1.2." -> "...1.1.3."
$d.Content.Find.Execute("1.1.2.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "1.1.3.", 2)

# 2. Append an explanatory sentence to the end of the AnimTool paragraph.
$d.Content.Find.Execute("This is a purely cosmetic issue and is fixed by changing properties.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, `
                         "This is a purely cosmetic issue and is fixed by changing properties. This issue has not been fixed yet because critical issues have been found that require immediate fixes, as well as the development of version 2.0 (renamed from version 1.2 due to a 0% API compatibility rate)", `
                         2)

# --- Header updates --------------------------------------------------------

$header = $d.Sections(1).Headers(1)

# 3. Bump the patch version in the header title:
#    "Lightning Known Issues for 1.1.2" -> "...1.1.3"
$header.Range.Find.Execute("Lightning Known Issues for 1.1.2", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "Lightning Known Issues for 1.1.3", 2)

# 4. Update the date shown in the header.
$header.Range.Find.Execute("November 5, 2022", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "November 8, 2022", 2)
